# Criando Conexao ao Banco de Dados
# Replace the two sample/placeholder rows in the "Clientes" sheet with a
# single record driven from the (new) database connection: only the
# person's name is written to column A of row 2; the old CPF column and
# the second sample row are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old placeholder data (A2:C3 -> "1/Gizelli/877.872.112-74"
# twice) so no stale cells/shared strings survive.
$ws.Range("A2:C3").ClearContents()

# Write the single new record coming from the database.
$ws.Range("A2").Value = "henderson"
